$wb = $excel.ActiveWorkbook

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsTimeLog = $wb.Worksheets.Item("time_log")

# --- time_log: add two new logged entries (rows 3 & 4) ---
$wsTimeLog.Range("A3").Value = 45057
$wsTimeLog.Range("B3").Value = "sample ID investigation"
$wsTimeLog.Range("C3").Value = 0.52

$wsTimeLog.Range("A4").Value = 45061
$wsTimeLog.Range("B4").Value = "Finalize implementation of fool-proof method to differentiate baseline from treated recording"
$wsTimeLog.Range("C4").Value = 1.78

# Match the date number format used by the existing date column (A2)
$wsTimeLog.Range("A2").Copy() | Out-Null
$wsTimeLog.Range("A3").PasteSpecial(-4122) | Out-Null
$wsTimeLog.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- update last-selected cell on each sheet ---
$wsSheet1.Range("A7").Select() | Out-Null
$wsTimeLog.Range("B5").Select() | Out-Null
